# Update Pre_Mean (F), Post_Mean (G), and Change (H) values for the
# dq30_pct_$ KPI rows (57-111) in the Branch Summary sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F57").Value = 0.8664842074263248
$ws.Range("G57").Value = 0.8617165690482793
$ws.Range("H57").Value = -0.004767638378045502
$ws.Range("F58").Value = 0.8686697444042524
$ws.Range("G58").Value = 0.869913188008732
$ws.Range("H58").Value = 0.001243443604479588
$ws.Range("F59").Value = 0.8269244736885339
$ws.Range("G59").Value = 0.8275369401549431
$ws.Range("H59").Value = 0.0006124664664091117
$ws.Range("F60").Value = 0.8447601001738534
$ws.Range("G60").Value = 0.8495607300163122
$ws.Range("H60").Value = 0.00480062984245877
$ws.Range("F61").Value = 0.8596820703230861
$ws.Range("G61").Value = 0.8755587191527726
$ws.Range("H61").Value = 0.01587664882968653
$ws.Range("F62").Value = 0.6947638223458864
$ws.Range("G62").Value = 0.664657891413697
$ws.Range("H62").Value = -0.03010593093218938
$ws.Range("F63").Value = 0.6910166407791919
$ws.Range("G63").Value = 0.7138625544225599
$ws.Range("H63").Value = 0.02284591364336797
$ws.Range("F64").Value = 0.7004534185338477
$ws.Range("G64").Value = 0.709376896842749
$ws.Range("H64").Value = 0.008923478308901278
$ws.Range("F65").Value = 0.7209550597876914
$ws.Range("G65").Value = 0.7158874831966042
$ws.Range("H65").Value = -0.005067576591087142
$ws.Range("F66").Value = 0.7126446100962088
$ws.Range("G66").Value = 0.7215562931873156
$ws.Range("H66").Value = 0.008911683091106792
$ws.Range("F67").Value = 0.7284411266196106
$ws.Range("G67").Value = 0.7374191907631018
$ws.Range("H67").Value = 0.008978064143491227
$ws.Range("F68").Value = 0.626581551583688
$ws.Range("G68").Value = 0.6551327133162764
$ws.Range("H68").Value = 0.02855116173258843
$ws.Range("F69").Value = 0.687309499515669
$ws.Range("G69").Value = 0.6281373368003978
$ws.Range("H69").Value = -0.05917216271527115
$ws.Range("F70").Value = 0.7734732755666378
$ws.Range("G70").Value = 0.7786272560307739
$ws.Range("H70").Value = 0.005153980464136065
$ws.Range("F71").Value = 0.8204474026564186
$ws.Range("G71").Value = 0.8421200180693439
$ws.Range("H71").Value = 0.02167261541292531
$ws.Range("F72").Value = 0.855026202099714
$ws.Range("G72").Value = 0.8734011922887865
$ws.Range("H72").Value = 0.01837499018907252
$ws.Range("F73").Value = 0.7788243738790469
$ws.Range("G73").Value = 0.7695620292236202
$ws.Range("H73").Value = -0.009262344655426746
$ws.Range("F74").Value = 0.7363621715291051
$ws.Range("G74").Value = 0.7588191354930586
$ws.Range("H74").Value = 0.02245696396395347
$ws.Range("F75").Value = 0.7593254759047275
$ws.Range("G75").Value = 0.7539219139167801
$ws.Range("H75").Value = -0.005403561987947447
$ws.Range("F76").Value = 0.6799881371319128
$ws.Range("G76").Value = 0.6874327865623868
$ws.Range("H76").Value = 0.007444649430474071
$ws.Range("F77").Value = 0.7054564693230709
$ws.Range("G77").Value = 0.7103585252747668
$ws.Range("H77").Value = 0.004902055951695905
$ws.Range("F78").Value = 0.7493116611707433
$ws.Range("G78").Value = 0.7587768963873039
$ws.Range("H78").Value = 0.009465235216560619
$ws.Range("F79").Value = 0.7063570713977596
$ws.Range("G79").Value = 0.721097304871533
$ws.Range("H79").Value = 0.0147402334737734
$ws.Range("F80").Value = 0.7335447817389327
$ws.Range("G80").Value = 0.7362218319216092
$ws.Range("H80").Value = 0.002677050182676499
$ws.Range("F81").Value = 0.6780682516659573
$ws.Range("G81").Value = 0.6830674077295056
$ws.Range("H81").Value = 0.004999156063548282
$ws.Range("F82").Value = 0.7142788750429839
$ws.Range("G82").Value = 0.741677946464538
$ws.Range("H82").Value = 0.02739907142155407
$ws.Range("F83").Value = 0.6498112635082486
$ws.Range("G83").Value = 0.6907347603567985
$ws.Range("H83").Value = 0.0409234968485499
$ws.Range("F84").Value = 0.6552365365958429
$ws.Range("G84").Value = 0.5729141456667907
$ws.Range("H84").Value = -0.08232239092905225
$ws.Range("F85").Value = 0.6792110123817749
$ws.Range("G85").Value = 0.6947787375161391
$ws.Range("H85").Value = 0.01556772513436411
$ws.Range("F86").Value = 0.7226945554241307
$ws.Range("G86").Value = 0.7352180953233294
$ws.Range("H86").Value = 0.01252353989919863
$ws.Range("F87").Value = 0.6594240389935497
$ws.Range("G87").Value = 0.6621442033231031
$ws.Range("H87").Value = 0.002720164329553354
$ws.Range("F88").Value = 0.6152344366478365
$ws.Range("G88").Value = 0.6511950853002084
$ws.Range("H88").Value = 0.03596064865237192
$ws.Range("F89").Value = 0.6833995400163368
$ws.Range("G89").Value = 0.6494359536503761
$ws.Range("H89").Value = -0.03396358636596075
$ws.Range("F90").Value = 0.6796962366538458
$ws.Range("G90").Value = 0.6729059695477085
$ws.Range("H90").Value = -0.006790267106137282
$ws.Range("F91").Value = 0.712346497114461
$ws.Range("G91").Value = 0.6814733687353028
$ws.Range("H91").Value = -0.03087312837915812
$ws.Range("F92").Value = 0.7355023839675536
$ws.Range("G92").Value = 0.7522899250826447
$ws.Range("H92").Value = 0.01678754111509106
$ws.Range("F93").Value = 0.5798335294544306
$ws.Range("G93").Value = 0.625196788915936
$ws.Range("H93").Value = 0.04536325946150543
$ws.Range("F94").Value = 0.7235019449790707
$ws.Range("G94").Value = 0.7070301315247397
$ws.Range("H94").Value = -0.01647181345433102
$ws.Range("F95").Value = 0.6807681776666373
$ws.Range("G95").Value = 0.716763973072162
$ws.Range("H95").Value = 0.03599579540552467
$ws.Range("F96").Value = 0.7136705597717454
$ws.Range("G96").Value = 0.7281785519211766
$ws.Range("H96").Value = 0.01450799214943121
$ws.Range("F97").Value = 0.6527559447536727
$ws.Range("G97").Value = 0.6288566170483629
$ws.Range("H97").Value = -0.02389932770530978
$ws.Range("F98").Value = 0.7769792242155896
$ws.Range("G98").Value = 0.7979242066375748
$ws.Range("H98").Value = 0.02094498242198517
$ws.Range("F99").Value = 0.7292172815929315
$ws.Range("G99").Value = 0.7150436056494053
$ws.Range("H99").Value = -0.01417367594352625
$ws.Range("F100").Value = 0.7038639079168134
$ws.Range("G100").Value = 0.6940904805720286
$ws.Range("H100").Value = -0.009773427344784813
$ws.Range("F101").Value = 0.6566331047878629
$ws.Range("G101").Value = 0.6483516460183039
$ws.Range("H101").Value = -0.008281458769558991
$ws.Range("F102").Value = 0.63917438835754
$ws.Range("G102").Value = 0.6355054009612549
$ws.Range("H102").Value = -0.003668987396285184
$ws.Range("F103").Value = 0.6769020360891681
$ws.Range("G103").Value = 0.6859448605708783
$ws.Range("H103").Value = 0.009042824481710232
$ws.Range("F104").Value = 0.7175202438336
$ws.Range("G104").Value = 0.7018538432857423
$ws.Range("H104").Value = -0.01566640054785773
$ws.Range("F105").Value = 0.7769535658727831
$ws.Range("G105").Value = 0.7508980978238089
$ws.Range("H105").Value = -0.02605546804897418
$ws.Range("F106").Value = 0.8030614300118006
$ws.Range("G106").Value = 0.7952301125909335
$ws.Range("H106").Value = -0.007831317420867179
$ws.Range("F107").Value = 0.8224485590239767
$ws.Range("G107").Value = 0.8794387555907924
$ws.Range("H107").Value = 0.0569901965668157
$ws.Range("F108").Value = 0.778168888643356
$ws.Range("G108").Value = 0.768064392944534
$ws.Range("H108").Value = -0.01010449569882199
$ws.Range("F109").Value = 0.8442698905143446
$ws.Range("G109").Value = 0.8500309356176029
$ws.Range("H109").Value = 0.005761045103258255
$ws.Range("F110").Value = 0.7635096610963941
$ws.Range("G110").Value = 0.7616381780663601
$ws.Range("H110").Value = -0.001871483030033927
$ws.Range("F111").Value = 0.7069551338827909
$ws.Range("G111").Value = 0.713758786297309
$ws.Range("H111").Value = 0.006803652414518191
